$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheet 1) — "想去人数" (F column) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 789
$ws1.Range("F5").Value = 144
$ws1.Range("F6").Value = 18
$ws1.Range("F7").Value = 167
$ws1.Range("F8").Value = 352
$ws1.Range("F9").Value = 463
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 11918
$ws1.Range("F13").Value = 5430

# Sheet "演出" (Worksheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 108

# Sheet "全部类型" (Worksheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 789
$ws4.Range("F4").Value = 108
$ws4.Range("F7").Value = 144
$ws4.Range("F8").Value = 18
$ws4.Range("F9").Value = 167
$ws4.Range("F10").Value = 352
$ws4.Range("F11").Value = 463
$ws4.Range("F13").Value = 144
$ws4.Range("F14").Value = 11918
$ws4.Range("F16").Value = 5430
